$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.06045809657875509
$ws.Cells.Item(2, 3).Value = 0.04715535591891339
$ws.Cells.Item(2, 4).Value = 0.02048158440005299
$ws.Cells.Item(2, 5).Value = 0.01686187421884393
$ws.Cells.Item(2, 6).Value = 0.03134567345318518
$ws.Cells.Item(3, 2).Value = 0.05401023734024143
$ws.Cells.Item(3, 3).Value = 0.04664377912446473
$ws.Cells.Item(3, 4).Value = 0.02187586734642404
$ws.Cells.Item(3, 5).Value = 0.01715629985382717
$ws.Cells.Item(3, 6).Value = 0.02956491716692805
$ws.Cells.Item(4, 2).Value = 0.0470616199553073
$ws.Cells.Item(4, 3).Value = 0.04567511983984629
$ws.Cells.Item(4, 4).Value = 0.02333571294534873
$ws.Cells.Item(4, 5).Value = 0.02089269369606267
$ws.Cells.Item(4, 6).Value = 0.03038370708262419
$ws.Cells.Item(5, 2).Value = 0.07308407875674522
$ws.Cells.Item(5, 3).Value = 0.0505855742851702
$ws.Cells.Item(5, 4).Value = 0.02866222179820695
$ws.Cells.Item(5, 5).Value = 0.03114140303348209
$ws.Cells.Item(5, 6).Value = 0.04009766902101392
$ws.Cells.Item(6, 2).Value = 0.07082668202974202
$ws.Cells.Item(6, 3).Value = 0.07054062889036722
$ws.Cells.Item(6, 4).Value = 0.03237154201126941
$ws.Cells.Item(6, 5).Value = 0.03370693571387755
$ws.Cells.Item(6, 6).Value = 0.04606262693942028
$ws.Cells.Item(7, 2).Value = 0.06169330846661639
$ws.Cells.Item(7, 3).Value = 0.08360431731307907
$ws.Cells.Item(7, 4).Value = 0.03022029982838398
$ws.Cells.Item(7, 5).Value = 0.02707951426573501
$ws.Cells.Item(7, 6).Value = 0.04999358501087372
$ws.Cells.Item(8, 2).Value = 0.06629066368454885
$ws.Cells.Item(8, 3).Value = 0.06135416057263635
$ws.Cells.Item(8, 4).Value = 0.02368861187480092
$ws.Cells.Item(8, 5).Value = 0.02170854605108558
$ws.Cells.Item(8, 6).Value = 0.03971103758967622
$ws.Cells.Item(9, 2).Value = 0.05216648029813976
$ws.Cells.Item(9, 3).Value = 0.05136802932237945
$ws.Cells.Item(9, 4).Value = 0.02039295078470096
$ws.Cells.Item(9, 5).Value = 0.01638579118884057
$ws.Cells.Item(9, 6).Value = 0.02933722974991203
$ws.Cells.Item(10, 2).Value = 0.06510615745948943
$ws.Cells.Item(10, 3).Value = 0.04680887516166399
$ws.Cells.Item(10, 4).Value = 0.02033912943808381
$ws.Cells.Item(10, 5).Value = 0.0134631626092732
$ws.Cells.Item(10, 6).Value = 0.03091019378675437
$ws.Cells.Item(11, 2).Value = 0.07545211636098324
$ws.Cells.Item(11, 3).Value = 0.04703330337588309
$ws.Cells.Item(11, 4).Value = 0.0219230065764918
$ws.Cells.Item(11, 5).Value = 0.01073546748923419
$ws.Cells.Item(11, 6).Value = 0.02903007714123914
$ws.Cells.Item(12, 2).Value = 0.07410596051833308
$ws.Cells.Item(12, 3).Value = 0.04608447413816349
$ws.Cells.Item(12, 4).Value = 0.02401308390205573
$ws.Cells.Item(12, 5).Value = 0.01208761325420784
$ws.Cells.Item(12, 6).Value = 0.03034085192087571
$ws.Cells.Item(13, 2).Value = 0.07702269859463579
$ws.Cells.Item(13, 3).Value = 0.04663826203447504
$ws.Cells.Item(13, 4).Value = 0.02473300078266149
$ws.Cells.Item(13, 5).Value = 0.01433096314591077
$ws.Cells.Item(13, 6).Value = 0.02993318911252935
$ws.Cells.Item(14, 2).Value = 0.07295338397274957
$ws.Cells.Item(14, 3).Value = 0.04715978001768463
$ws.Cells.Item(14, 4).Value = 0.02632171017002597
$ws.Cells.Item(14, 5).Value = 0.01507888315153952
$ws.Cells.Item(14, 6).Value = 0.03053072992730379
$ws.Cells.Item(15, 2).Value = 0.07584408747861211
$ws.Cells.Item(15, 3).Value = 0.04682685968381877
$ws.Cells.Item(15, 4).Value = 0.0259957302141631
$ws.Cells.Item(15, 5).Value = 0.01524926703987872
$ws.Cells.Item(15, 6).Value = 0.03050423236089556
$ws.Cells.Item(16, 2).Value = 0.07627834004907812
$ws.Cells.Item(16, 3).Value = 0.04690048807944901
$ws.Cells.Item(16, 4).Value = 0.02625434332536366
$ws.Cells.Item(16, 5).Value = 0.01568058590965807
$ws.Cells.Item(16, 6).Value = 0.03075172629138184
$ws.Cells.Item(17, 2).Value = 0.0757725526778312
$ws.Cells.Item(17, 3).Value = 0.04680293220379636
$ws.Cells.Item(17, 4).Value = 0.02649785703272504
$ws.Cells.Item(17, 5).Value = 0.01591807867429529
$ws.Cells.Item(17, 6).Value = 0.0309157460626473
$ws.Cells.Item(18, 2).Value = 0.07579929206477883
$ws.Cells.Item(18, 3).Value = 0.04695387455426898
$ws.Cells.Item(18, 4).Value = 0.0265071244776255
$ws.Cells.Item(18, 5).Value = 0.01602796623784324
$ws.Cells.Item(18, 6).Value = 0.03091670582609169
